$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.22"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.99"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.334"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05962"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.393"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8110"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9638"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1425"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07382"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03385"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03051"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09406"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.993"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001599"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04815"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005912"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006254"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005138"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009836"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.00009703"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.743"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.185"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03911"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006442"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1071"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005831"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005316"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8503"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.03601"
